# Major edit for cost reduction:
# Inserts three new parameter rows into the "Sheet1" parameter table:
#   - carrier_y2 (new row before carrier_z)
#   - carrier_z2 (new row right after carrier_z)
#   - frame_x    (new row between frame_total_z and frame_y)
# and rewires frame_border's formula to use the new frame_x parameter.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- 1) Insert "carrier_y2" row above row 25 (old carrier_z row) ---
$ws.Rows.Item(25).Insert()
$ws.Range("A25").Value = "carrier_y2"
$ws.Range("B25").Formula = "=B64-B2*2"
$ws.Range("D25").Value = "new"

# --- 2) Insert "carrier_z2" row below carrier_z (now row 26) ---
$ws.Rows.Item(27).Insert()
$ws.Range("A27").Value = "carrier_z2"
$ws.Range("B27").Value = 25
$ws.Range("D27").Value = "new"

# --- 3) Insert "frame_x" row between frame_total_z and frame_y ---
# (before these inserts frame_total_z/frame_y were rows 60/61; after the two
#  row inserts above they are rows 62/63, so the new row goes in at 63)
$ws.Rows.Item(63).Insert()
$ws.Range("A63").Value = "frame_x"
$ws.Range("B63").Formula = "=B23+0.2+B2*2"
$ws.Range("D63").Value = "new"

# --- 4) Rewrite frame_border's formula (now row 66) to depend on frame_x ---
$ws.Range("B66").Formula = "=(B63-(B48+B9*2))/2"

# --- 5) Match the saved view state from the edit ---
$ws.Range("B67").Select()
$ws.Application.ActiveWindow.ScrollRow = 43
